# Generate Report for Handback
# Marks the two files in each locale sheet (zh-cn, de-de) as handed back:
#  - Overview sheet status columns -> "Handed back: in sync with en-US"
#  - Per-locale sheets: Status column, Latest Target File, Latest Handback
#    File and Latest Handback DateTime columns are filled in, and the
#    Latest Target File cell is turned into a hyperlink to the source .md
#    file (matching column A's hyperlink).

$wb = $excel.ActiveWorkbook

$statusOld = "In Translation"
$statusNew = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: columns E (zh-cn) and F (de-de) show the same status
# text for each of the two source files (rows 2-3).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusNew
$overview.Range("F2").Value = $statusNew
$overview.Range("E3").Value = $statusNew
$overview.Range("F3").Value = $statusNew
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# Per-locale sheets: "zh-cn" and "de-de" both have the same table
# layout (column letters) but different handback file names / dates.
# ---------------------------------------------------------------------
$locales = @(
    @{
        SheetName   = "zh-cn"
        Row2Target  = "429986db-8279-4169-8796-64c3284c0028.c58bbb4d058f3d4bb8d188a8b51a04c58bd8a2cd.zh-cn.xlf"
        Row3Target  = "4a03e220-4bd4-40b1-8232-87ee5231bbdb.0a3db43e27db1ad45976085f96883bde4a05ad20.zh-cn.xlf"
        HandbackDateTime = "2016-09-04 14:27:52"
    },
    @{
        SheetName   = "de-de"
        Row2Target  = "429986db-8279-4169-8796-64c3284c0028.c58bbb4d058f3d4bb8d188a8b51a04c58bd8a2cd.de-de.xlf"
        Row3Target  = "4a03e220-4bd4-40b1-8232-87ee5231bbdb.0a3db43e27db1ad45976085f96883bde4a05ad20.de-de.xlf"
        HandbackDateTime = "2016-09-04 14:27:59"
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.SheetName)

    # Column widths widened to fit the new hyperlink / longer values.
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40

    # Existing hyperlinks on column A (row 2 -> file 1, row 3 -> file 2)
    # give us the target URL + display text to mirror onto column I.
    # (Iterate the worksheet-level Hyperlinks collection -- indexing via
    # Range.Hyperlinks does not reliably surface Address/TextToDisplay.)
    $linkMap = @{}
    foreach ($h in $ws.Hyperlinks) {
        $linkMap[$h.Range.Address()] = @{ Address = $h.Address; Text = $h.TextToDisplay }
    }
    $linkA2 = $linkMap['$A$2']
    $linkA3 = $linkMap['$A$3']
    $addrA2 = $linkA2.Address
    $textA2 = $linkA2.Text
    $addrA3 = $linkA3.Address
    $textA3 = $linkA3.Text

    # Row 2 (first source file)
    $ws.Range("C2").Value = $statusNew
    $ws.Range("I2").Value = $textA2
    $ws.Hyperlinks.Add($ws.Range("I2"), $addrA2, "", "", $textA2)
    $ws.Range("J2").Value = $locale.Row2Target
    $ws.Range("K2").Value = $locale.HandbackDateTime

    # Row 3 (second source file)
    $ws.Range("C3").Value = $statusNew
    $ws.Range("I3").Value = $textA3
    $ws.Hyperlinks.Add($ws.Range("I3"), $addrA3, "", "", $textA3)
    $ws.Range("J3").Value = $locale.Row3Target
    $ws.Range("K3").Value = $locale.HandbackDateTime
}
